$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New normalized damage values for B2:B37
$values = @(
    28.94698445859289,
    39.558216993653311,
    34.762463242137294,
    34.006513930276817,
    29.415956505802455,
    31.7640372787662,
    34.257855442714273,
    30.934095646603833,
    30.693298240022809,
    31.441756861747034,
    32.745543391442453,
    31.468762161095309,
    29.641657203803394,
    30.706180483338958,
    31.088782082708889,
    32.101628072176005,
    32.628603853396868,
    32.058498703461865,
    31.910453727297021,
    29.868455489280382,
    30.795206964065375,
    28.209242410455527,
    29.26068341235661,
    33.429666999767832,
    31.014346202994169,
    29.556480501294853,
    34.332259127991954,
    33.247400842597152,
    31.368414518660348,
    30.026679125872779,
    31.034163781265807,
    37.729126660417897,
    35.437380208660834,
    30.01216891788005,
    34.015224646399176,
    100
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Swap the "Estado" labels: row 2 goes from "Daño" to " - ", row 37 goes from " - " to "Daño"
$ws.Cells.Item(2, 3).Value = " - "
$ws.Cells.Item(37, 3).Value = "Daño"
